# Update countries & provincias Spain
# Refresh the "last updated" timestamp shown in A1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 12:57"

# --- Row 4 (Estados Unidos): updated case counts ---
$ws.Range("B4").Value = 5251997
$ws.Range("C4").Value = 551
$ws.Range("D4").Value = 2716593
$ws.Range("E4").Value = 2369203
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 166201

# --- Row 14 (Iran): updated case counts ---
$ws.Range("B14").Value = 331189
$ws.Range("C14").Value = 2345
$ws.Range("D14").Value = 288620
$ws.Range("E14").Value = 23769
$ws.Range("G14").Value = 184
$ws.Range("H14").Value = 18800

# --- Row 42 (Bielorrusia): updated case counts ---
$ws.Range("B42").Value = 69005
$ws.Range("C42").Value = 58
$ws.Range("D42").Value = 65219
$ws.Range("E42").Value = 3194
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 592

# --- Rows 43/44: Rumania overtakes Emiratos Arabes Unidos, so the two rows swap countries ---
$ws.Range("A43").Value = "Rumania"
$ws.Range("B43").Value = 63762
$ws.Range("C43").Value = 1215
$ws.Range("D43").Value = 30585
$ws.Range("E43").Value = 30413
$ws.Range("G43").Value = 35
$ws.Range("H43").Value = 2764

$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("B44").Value = 62966
$ws.Range("C44").Value = 262
$ws.Range("D44").Value = 56961
$ws.Range("E44").Value = 5647
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 358

# --- Row 58 (Suiza): updated case counts ---
$ws.Range("B58").Value = 36895
$ws.Range("C58").Value = 187
$ws.Range("E58").Value = 2505
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 1990

# --- Rows 68/69/70: Nepal overtakes Costa Rica and Etiopia, rows re-labelled accordingly ---
$ws.Range("A68").Value = "Nepal"
$ws.Range("B68").Value = 23948
$ws.Range("C68").Value = 638
$ws.Range("D68").Value = 16664
$ws.Range("E68").Value = 7201
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 83

$ws.Range("A69").Value = "Costa Rica"
$ws.Range("B69").Value = 23872
$ws.Range("D69").Value = 7823
$ws.Range("E69").Value = 15805
$ws.Range("H69").Value = 244

$ws.Range("A70").Value = "Etiopia"
$ws.Range("B70").Value = 23591
$ws.Range("D70").Value = 10411
$ws.Range("E70").Value = 12760
$ws.Range("H70").Value = 420

# --- Row 82 (Madagascar): updated case counts ---
$ws.Range("B82").Value = 13317
$ws.Range("C82").Value = 115
$ws.Range("D82").Value = 11276
$ws.Range("E82").Value = 1889
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 152

# --- Row 85 (Senegal): updated case counts ---
$ws.Range("B85").Value = 11380
$ws.Range("C85").Value = 68
$ws.Range("D85").Value = 7449
$ws.Range("E85").Value = 3693
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 238

# --- Row 87 (Consejo Danes para los Refugiados): updated case counts ---
$ws.Range("B87").Value = 9499
$ws.Range("C87").Value = 10
$ws.Range("D87").Value = 8375
$ws.Range("E87").Value = 899
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 225

# --- Row 105 (Maldivas): updated case counts ---
$ws.Range("E105").Value = 2302
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 20

# --- Row 122 (Sri Lanka): updated case counts ---
$ws.Range("B122").Value = 2875
$ws.Range("C122").Value = 4
$ws.Range("E122").Value = 242

# --- Row 141 (Uganda): updated case counts ---
$ws.Range("B141").Value = 1313
$ws.Range("C141").Value = 16
$ws.Range("D141").Value = 1138
$ws.Range("E141").Value = 166

# --- Row 151 (Malta): updated case counts ---
$ws.Range("B151").Value = 1141
$ws.Range("C151").Value = 29
$ws.Range("D151").Value = 692
$ws.Range("E151").Value = 440

# --- Rows 213/214: Montserrat and Islas Malvinas swap (tie-break reorder) ---
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
